$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 50 (Modelo MCPS, d=1)
$ws.Range("C50").Value = 9.870731880388252
$ws.Range("D50").Value = 1.528800641636963 / 1000000000

# Row 51 (Modelo MCPS, d=2)
$ws.Range("C51").Value = 4.45838598591608
$ws.Range("D51").Value = 0.0001969007947060231

# Row 52 (Modelo MCPS, d=3)
$ws.Range("C52").Value = 5.762959678840502
$ws.Range("D52").Value = 8.49495688237667 / 1000000

# Row 53 (Modelo MCPS, d=4)
$ws.Range("C53").Value = 6.037894737575589
$ws.Range("D53").Value = 4.461552906676047 / 1000000

# Row 54 (Modelo MCPS, d=5)
$ws.Range("C54").Value = 5.975061803818217
$ws.Range("D54").Value = 5.16519300819418 / 1000000

# Row 55 (Modelo MCPS, d=6)
$ws.Range("C55").Value = 8.033005460928365
$ws.Range("D55").Value = 5.512877110369629 / 100000000

# Row 56 (Modelo MCPS, d=7)
$ws.Range("C56").Value = 5.786140801158132
$ws.Range("D56").Value = 8.043497842891156 / 1000000

# Row 57 (Modelo MCPS, d=10)
$ws.Range("C57").Value = 7.226498109201343
$ws.Range("D57").Value = 3.058914341824703 / 10000000
